# Training/validation log sheet: update the DATE column (column B) for all
# data rows (rows 2-80) from "2024-09-20" to "2024-10-01".
#
# The DATE column stores its values as literal text (e.g. "2024-09-20"),
# not as real Excel date serials. Writing a date-shaped string straight
# into a General-formatted cell would normally get auto-converted into a
# date serial by Excel's input parser, so the column is first marked as
# Text ("@") to force the replacement to stick as plain text, matching
# the original data shape.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$dateRange = $ws.Range("B2:B80")
$dateRange.NumberFormat = "@"

for ($row = 2; $row -le 80; $row++) {
    $cell = $ws.Cells.Item($row, 2)
    if ($cell.Text -eq "2024-09-20") {
        $cell.Value = "2024-10-01"
    }
}
